$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.516.22"
$ws.Range("E2").Value = "  +2.86%  "

$ws.Range("D3").Value = "2.414.26"
$ws.Range("E3").Value = "  +8.92%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'324.26"
$ws.Range("E5").Value = "  +12.86%  "

$ws.Range("D6").Value = "'104.61"
$ws.Range("E6").Value = "  -3.93%  "

$ws.Range("E7").Value = "  +4.29%  "

$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("E9").Value = "  +12.36%  "

$ws.Range("D10").Value = "'42.25"
$ws.Range("E10").Value = "  -1.35%  "

$ws.Range("D11").Value = "'0.0950"
$ws.Range("E11").Value = "  +5.28%  "

$ws.Range("D12").Value = "'8.66"
$ws.Range("E12").Value = "  +2.13%  "

$ws.Range("D13").Value = "'1.04"
$ws.Range("E13").Value = "  +4.53%  "

$ws.Range("D14").Value = "'17.27"
$ws.Range("E14").Value = "  +17.06%  "

$ws.Range("E15").Value = "  +2.91%  "

$ws.Range("D16").Value = "2.780.24"
$ws.Range("E16").Value = "  +8.93%  "

$ws.Range("D17").Value = "2.418.32"
$ws.Range("E17").Value = "  +8.36%  "

$ws.Range("D18").Value = "43.555.88"
$ws.Range("E18").Value = "  +3.13%  "

$ws.Range("D19").Value = "'0.0000110"
$ws.Range("E19").Value = "  +6.47%  "

$ws.Range("D20").Value = "'7.46"
$ws.Range("E20").Value = "  +5.60%  "

$ws.Range("D21").Value = "'75.54"
$ws.Range("E21").Value = "  +4.27%  "

$ws.Range("E22").Value = "  +4.64%  "

$ws.Range("D23").Value = "'260.31"
$ws.Range("E23").Value = "  +14.12%  "

$ws.Range("E24").Value = "  +5.05%  "

$ws.Range("D25").Value = "'9.67"
$ws.Range("E25").Value = "  +9.83%  "

$ws.Range("D26").Value = "'11.96"
$ws.Range("E26").Value = "  +6.12%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").Value = "'22.84"
$ws.Range("E29").Value = "  +10.78%  "

$ws.Range("D30").Value = "'179.75"
$ws.Range("E30").Value = "  +4.22%  "

$ws.Range("D31").Value = "'2.23"
$ws.Range("E31").Value = "  +1.84%  "

$ws.Range("D32").Value = "'38.22"
$ws.Range("E32").Value = "  +4.44%  "

$ws.Range("E33").Value = "  +2.68%  "

$ws.Range("D34").Value = "'0.0935"
$ws.Range("E34").Value = "  +8.26%  "

$ws.Range("E35").Value = "  +7.65%  "

$ws.Range("E36").Value = "  +6.67%  "

$ws.Range("D37").Value = "'4.89"
$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("E38").Value = "  +2.64%  "

$ws.Range("D39").Value = "'3.99"
$ws.Range("E39").Value = "  -2.37%  "

$ws.Range("D40").Value = "'0.106"
$ws.Range("E40").Value = "  +4.29%  "

$ws.Range("D41").Value = "'2.89"
$ws.Range("E41").Value = "  +23.13%  "

$ws.Range("D42").Value = "'1.62"
$ws.Range("E42").Value = "  +26.96%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'126.19"
$ws.Range("E43").Value = "  +26.14%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.234"
$ws.Range("E44").Value = "  +3.58%  "

$ws.Range("D45").Value = "'69.85"
$ws.Range("E45").Value = "  -4.05%  "

$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "'12.61"
$ws.Range("E47").Value = "  +4.29%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'9.57"
$ws.Range("E48").Value = "  +14.66%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'5.69"
$ws.Range("E49").Value = "  +7.91%  "

$ws.Range("E50").Value = "  +5.73%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.606.57"
$ws.Range("E51").Value = "  +14.60%  "
